{"js": "// Load all paragraphs in the body so we can locate the \"Xem th\u00f4ng b\u00e1o\" item.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that holds the \"Xem th\u00f4ng b\u00e1o\" bullet and insert a new\n// bullet right after it (inherits the same ListParagraph / numbering style).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Xem th\u00f4ng b\u00e1o\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  target.insertParagraph(\"Th\u00f4ng b\u00e1o realtime\", Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Remove the leftover \"_GoBack\" bookmark (Word regenerates this automatically\n// while editing; it is not meant to be kept in the saved document).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Find the bullet paragraph that reads \"Xem th\u00f4ng b\u00e1o\" and add a new bullet\n# right after it (new paragraph inherits the same ListParagraph / numbering\n# formatting from the paragraph it was split off from).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Xem th\u00f4ng b\u00e1o\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $index = $target.Index\n    $target.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs($index + 1)\n    $newPara.Range.Text = \"Th\u00f4ng b\u00e1o realtime\"\n}\n\n# Remove the leftover \"_GoBack\" bookmark (Word regenerates this automatically\n# while editing; it is not meant to be kept in the saved document).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
